$p = $ppt.ActivePresentation
$s = $p.Slides.Add(6, 12)
$sh = $s.Shapes.AddShape(1, 174, 72, 72, 72)
$sh.Fill.ForeColor.ObjectThemeColor = 5
$sh.Fill.ForeColor.RGB = 0xF1E5DB
Write-Output ("type=" + $sh.Fill.ForeColor.Type + " rgb=" + $sh.Fill.ForeColor.RGB + " theme=" + $sh.Fill.ForeColor.ObjectThemeColor)
